$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D12", "D15", "D16", "D18", "D19", "D20", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.320.86"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.931.19"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "251.76"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").Value = "0.7127"
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "0.3261"
$ws.Range("E8").Value = "  -0.47%  "
$ws.Range("D9").Value = "27.52"
$ws.Range("E9").Value = "  +3.29%  "
$ws.Range("D10").Value = "0.07203"
$ws.Range("E10").Value = "  +5.67%  "
$ws.Range("D11").Value = "0.7990"
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").Value = "0.08090"
$ws.Range("E12").Value = "  +1.70%  "
$ws.Range("D13").Value = "1.927.61"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").Value = "94.69"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "14.83"
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("D17").Value = "30.294.62"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "251.41"
$ws.Range("E18").Value = "  -3.89%  "
$ws.Range("D19").Value = "0.000008119"
$ws.Range("E19").Value = "  +2.25%  "
$ws.Range("D20").Value = "5.784"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").Value = "2.179.84"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Value = "6.922"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").Value = "9.717"
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("D26").Value = "164.74"
$ws.Range("E26").Value = "  +3.02%  "
$ws.Range("D27").Value = "19.20"
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("D28").Value = "2.315"
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("D29").Value = "0.1285"
$ws.Range("E29").Value = "  -4.65%  "
$ws.Range("D30").Value = "1.360"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").Value = "1.545"
$ws.Range("D32").Value = "4.431"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("D34").Value = "0.05210"
$ws.Range("E34").Value = "  +2.41%  "
$ws.Range("D35").Value = "1.265"
$ws.Range("E35").Value = "  +5.35%  "
$ws.Range("D36").Value = "0.7482"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").Value = "2.761"
$ws.Range("E37").Value = "  +1.29%  "
$ws.Range("D38").Value = "0.01961"
$ws.Range("E38").Value = "  +1.18%  "
$ws.Range("D39").Value = "2.796"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").Value = "78.88"
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("D42").Value = "0.4524"
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("D43").Value = "2.027"
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "0.8399"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("D46").Value = "101.83"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("D47").Value = "9.838"
$ws.Range("E47").Value = "  +1.13%  "
$ws.Range("D48").Value = "7.411"
$ws.Range("E48").Value = "  +1.76%  "
$ws.Range("D49").Value = "36.69"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").Value = "0.06088"
$ws.Range("E50").Value = "  +2.87%  "
$ws.Range("E51").Value = "  +1.43%  "
